# Daily attendance processing - 2025-12-27 19:28:20
# Reorder the "Recorded By" (column G) email/author lists for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$updates = @{
    2  = "servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, System, Amira.Sobhy@med.asu.edu.eg"
    3  = "majorelle.magdy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, System, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
    4  = "majorelle.magdy@med.asu.edu.eg, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
    5  = "Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
    6  = "majorelle.magdy@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, mennatulla.medhat@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, manar.montaser@med.asu.edu.eg"
    7  = "Fatmaelhady@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg"
    8  = "AbeerRagheb@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg"
    11 = "aml.awwad@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, aya.saeed@med.asu.edu.eg"
    12 = "dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg"
    13 = "amira.m.ibrahim@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg"
    15 = "mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
    17 = "mohamed.saleem@med.asu.edu.eg, esraa.sami@med.asu.edu.eg"
    25 = "menna-allah.gamil@med.asu.edu.eg, Noran.Mahmoud@med.asu.edu.eg"
    30 = "aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
